$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking strings in Price/Volume columns stay as text,
# matching the original inlineStr cell content instead of being coerced to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "67.064.31"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "3.494.34"
$ws.Range("E3").Value = "  -0.33%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "599.42"
$ws.Range("E5").Value = "  +0.60%  "
$ws.Range("D6").Value = "174.85"
$ws.Range("E6").Value = "  +2.72%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "0.586"
$ws.Range("E8").Value = "  -1.08%  "
$ws.Range("D9").Value = "0.130"
$ws.Range("E9").Value = "  -3.36%  "
$ws.Range("D10").Value = "7.15"
$ws.Range("E10").Value = "  -2.51%  "
$ws.Range("E11").Value = "  -1.47%  "
$ws.Range("D12").Value = "4.102.43"
$ws.Range("E12").Value = "  -0.22%  "
$ws.Range("D13").Value = "31.38"
$ws.Range("E13").Value = "  +10.40%  "
$ws.Range("E14").Value = "  +0.30%  "
$ws.Range("D15").Value = "67.056.19"
$ws.Range("E15").Value = "  +0.31%  "
$ws.Range("D16").Value = "0.0000178"
$ws.Range("E16").Value = "  -3.04%  "
$ws.Range("D17").Value = "3.492.45"
$ws.Range("E17").Value = "  -0.67%  "
$ws.Range("D18").Value = "6.28"
$ws.Range("E18").Value = "  -1.14%  "
$ws.Range("D19").Value = "14.49"
$ws.Range("E19").Value = "  +2.30%  "
$ws.Range("D20").Value = "393.23"
$ws.Range("E20").Value = "  -0.93%  "
$ws.Range("D21").Value = "7.98"
$ws.Range("E21").Value = "  -0.21%  "
$ws.Range("D22").Value = "73.36"
$ws.Range("E22").Value = "  -0.26%  "
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("E24").Value = "  -0.35%  "
$ws.Range("D25").Value = "5.69"
$ws.Range("E25").Value = "  -0.20%  "
$ws.Range("E26").Value = "  -2.77%  "
$ws.Range("D27").Value = "10.19"
$ws.Range("E27").Value = "  -0.61%  "
$ws.Range("E28").Value = "  -0.42%  "
$ws.Range("E29").Value = "  -0.49%  "
$ws.Range("D30").Value = "6.10"
$ws.Range("E30").Value = "  -3.43%  "
$ws.Range("E31").Value = "  -3.30%  "
$ws.Range("E32").Value = "  -0.58%  "
$ws.Range("E33").Value = "  -1.09%  "
$ws.Range("E34").Value = "  -0.53%  "
$ws.Range("D35").Value = "1.63"
$ws.Range("E35").Value = "  +0.98%  "
$ws.Range("D36").Value = "163.36"
$ws.Range("E36").Value = "  -0.25%  "
$ws.Range("E37").Value = "  -2.50%  "
$ws.Range("E38").Value = "  -0.38%  "
$ws.Range("D39").Value = "7.03"
$ws.Range("E39").Value = "  +2.24%  "
$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").Value = "4.66"
$ws.Range("E40").Value = "  -1.71%  "
$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").Value = "27.41"
$ws.Range("E41").Value = "  +0.90%  "
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").Value = "26.11"
$ws.Range("E42").Value = "  -1.73%  "
$ws.Range("B43").Value = "Hedera"
$ws.Range("C43").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D43").Value = "0.0729"
$ws.Range("E43").Value = "  -2.26%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "2.797.57"
$ws.Range("E44").Value = "  -1.21%  "
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").Value = "42.49"
$ws.Range("E45").Value = "  -0.79%  "
$ws.Range("B46").Value = "dogwifhat"
$ws.Range("C46").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D46").Value = "2.54"
$ws.Range("E46").Value = "  -4.44%  "
$ws.Range("D47").Value = "0.0302"
$ws.Range("E47").Value = "  -3.98%  "
$ws.Range("D48").Value = "337.18"
$ws.Range("E48").Value = "  -1.37%  "
$ws.Range("E49").Value = "  -2.19%  "
$ws.Range("D50").Value = "33.59"
$ws.Range("E50").Value = "  -0.90%  "
$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").Value = "6.44"
$ws.Range("E51").Value = "  -1.30%  "
